$d = $word.ActiveDocument

# Ordered list of (old, new) bullet/heading text replacements taken from the
# unified diff. Each entry is applied by locating the existing paragraph text
# with Find.Execute (no in-place Replace) and then assigning the new text to
# the matched Range directly. Doing the substitution this way - rather than
# via Find.Execute's own Replace argument - avoids Word's AutoCorrect engine
# 'smart quote' substitution mangling straight apostrophes in the replacement
# text (e.g. "company's") when going through Replace=wdReplaceAll.
#
# NOTE: "Political Research and Data Analysis" appears twice in the source
# document (RESEARCH DIRECTOR - PCCC section, and PROGRAMMER - Lake Research
# Partners section). Only the second occurrence is renamed by the diff, so it
# is handled separately below via direct paragraph indexing instead of a
# document-wide Find (which would otherwise hit the first occurrence too).
$replacements = @(
    @{ Old = '• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations'; New = '• Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions' }
    @{ Old = '• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics'; New = '• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics and demographic analysis' }
    @{ Old = '• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets'; New = '• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets' }
    @{ Old = '• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering'; New = '• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering' }
    @{ Old = '• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications'; New = '• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications' }
    @{ Old = '• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices'; New = '• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices' }
    @{ Old = '• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES'; New = '• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES' }
    @{ Old = '• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions'; New = '• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions' }
    @{ Old = '• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI'; New = '• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI' }
    @{ Old = '• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company''s distinguishing products'; New = '• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company''s distinguishing products' }
    @{ Old = '• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices'; New = '• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices' }
    @{ Old = '• Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies'; New = '• Architected and developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies' }
    @{ Old = '• Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu'; New = '• Collaborated with senior officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to enhance platform using Docker and Ubuntu' }
    @{ Old = '• Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions'; New = '• Conceived and developed predictive application using Python, Pandas, and Jupyter to forecast how crisis economies respond to different humanitarian interventions' }
    @{ Old = '• Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting'; New = '• Architected and developed RACSO, a comprehensive web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting' }
    @{ Old = '• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner'; New = '• Led RFP process and analyzed bids from 1,200 vendors before selecting optimal implementation partner' }
    @{ Old = '• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research'; New = '• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions' }
    @{ Old = '• Designed survey deployment system facilitating thousands of simultaneous phone surveys'; New = '• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs' }
    @{ Old = '• Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system'; New = '• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system' }
    @{ Old = '• Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill'; New = '• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill' }
    @{ Old = '• Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization'; New = '• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO' }
    @{ Old = '• Made all technology decisions and practices for massive multinational non-governmental organization'; New = '• Directed all technology decisions and practices for massive multinational non-governmental organization' }
    @{ Old = '• Wrote comprehensive frameworks for internal and external technology audits'; New = '• Developed comprehensive frameworks for internal and external technology audits' }
    @{ Old = '• Trained beneficiaries on spatial and Census data analysis for public health research'; New = '• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research' }
    @{ Old = '• Trained NGO staff in web development using Drupal, PHP, and MySQL'; New = '• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL' }
    @{ Old = '• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections'; New = '• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions' }
    @{ Old = 'Political Field Operations and Data Management'; New = 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns' }
    @{ Old = '• Administered all quantitative and qualitative research operations ensuring reporting accuracy'; New = '• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions' }
    @{ Old = '• Managed comprehensive survey fielding for multi-million dollar research firm'; New = '• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm' }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $found = $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $r.New
    } else {
        Write-Output "NOT FOUND: $($r.Old)"
    }
}

# The heading 'Political Research and Data Analysis' appears twice; only the
# second occurrence (PROGRAMMER - Lake Research Partners section) becomes
# 'Political Polling, Focus Groups and Demographic Analysis for Democratic
# Campaigns' per the diff. The first occurrence (RESEARCH DIRECTOR - PCCC)
# must stay untouched, so walk Paragraphs directly and only update the 2nd hit.
$targetOld = 'Political Research and Data Analysis'
$targetNew = 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns'
$matchIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetOld) {
        $matchIndex = $matchIndex + 1
        if ($matchIndex -eq 2) {
            $p.Range.Text = $targetNew
        }
    }
}
if ($matchIndex -ne 2) {
    Write-Output "WARNING: expected 2 occurrences of target heading, found $matchIndex"
}

Write-Output "Done."
